$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Questions Tracker")
$ws.Activate()

# Fill in the new row of data (row 48) that was previously blank
$ws.Range("B48").Value = "283. Move Zeroes"
$ws.Range("D48").Value = 1
$ws.Range("E48").Value = "19/11/2022"
$ws.Range("F48").Value = 1
$ws.Range("G48").Value = "Yes"

# Add the LeetCode hyperlink on C48 (also sets its display text/style)
$ws.Hyperlinks.Add($ws.Range("C48"), "https://leetcode.com/problems/move-zeroes/", "", "", "LeetCode")

# Move the active selection to H48, matching the saved view state
$ws.Range("H48").Select()

$wb.Save()
